$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "C:/Users/Asus/Desktop/21212/005\2504_image_005.jpg"
$ws.Range("D3").Value = "C:/Users/Asus/Desktop/21212/005\2505_image_005.jpg"
$ws.Range("D4").Value = "C:/Users/Asus/Desktop/21212/005\2506_image_005.jpg"
$ws.Range("D5").Value = "C:/Users/Asus/Desktop/21212/005\2507_image_005.jpg"
$ws.Range("D6").Value = "C:/Users/Asus/Desktop/21212/005\2508_image_005.jpg"
$ws.Range("D7").Value = "C:/Users/Asus/Desktop/21212/005\2509_image_005.jpg"
$ws.Range("D8").Value = "C:/Users/Asus/Desktop/21212/005\2510_image_005.jpg"
